$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$shape = $s.Shapes.Item(41)
$shape.Width = 68.1845
$shape.TextFrame.TextRange.Text = "qb:dataSet"
